# Adapt the column header formatting to the respective input file names:
#   *_old -> *_FV2304   (left-hand / "from" comparison columns)
#   *_new -> *_FV2310   (right-hand / "to" comparison columns)
# then promote the header + data range to a real Excel Table and freeze
# the header row so it stays visible while scrolling.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column headers in sheet order (A1:U1), old name -> new name.
$headerRenames = @(
    @{ old = "Segmentname_old";            new = "Segmentname_FV2304" },
    @{ old = "Segmentgruppe_old";          new = "Segmentgruppe_FV2304" },
    @{ old = "Segment_old";                new = "Segment_FV2304" },
    @{ old = "Datenelement_old";           new = "Datenelement_FV2304" },
    @{ old = "Segment ID_old";             new = "Segment ID_FV2304" },
    @{ old = "Code_old";                   new = "Code_FV2304" },
    @{ old = "Qualifier_old";              new = "Qualifier_FV2304" },
    @{ old = "Beschreibung_old";           new = "Beschreibung_FV2304" },
    @{ old = "Bedingungsausdruck_old";     new = "Bedingungsausdruck_FV2304" },
    @{ old = "Bedingung_old";              new = "Bedingung_FV2304" },
    @{ old = "diff";                       new = "diff" },
    @{ old = "Segmentname_new";            new = "Segmentname_FV2310" },
    @{ old = "Segmentgruppe_new";          new = "Segmentgruppe_FV2310" },
    @{ old = "Segment_new";                new = "Segment_FV2310" },
    @{ old = "Datenelement_new";           new = "Datenelement_FV2310" },
    @{ old = "Segment ID_new";             new = "Segment ID_FV2310" },
    @{ old = "Code_new";                   new = "Code_FV2310" },
    @{ old = "Qualifier_new";              new = "Qualifier_FV2310" },
    @{ old = "Beschreibung_new";           new = "Beschreibung_FV2310" },
    @{ old = "Bedingungsausdruck_new";     new = "Bedingungsausdruck_FV2310" },
    @{ old = "Bedingung_new";              new = "Bedingung_FV2310" }
)

for ($i = 0; $i -lt $headerRenames.Length; $i++) {
    $cell = $ws.Cells.Item(1, $i + 1)
    $cell.Value = $headerRenames[$i].new
}

# The header row (row 1) plus every data row make up the table range.
$tableRange = $ws.UsedRange

# Turn the range into a native Excel Table ("Table1") with an auto filter,
# mirroring the workbook's generated xl/tables/table1.xml.
$listObject = $ws.ListObjects.Add(1, $tableRange, $null, 1)
$listObject.Name = "Table1"
$listObject.TableStyle = ""

# Freeze the header row (row 1) so it stays pinned while scrolling.
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
